$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(21, 1100, 43.0975227355957,  4.08340311050415,  29.94179916381836, 57.77489471435547, 18476, "07"),
    @(22, 1200, 43.16032791137695, 4.144697666168213, 29.83244895935059, 57.38874816894531, 18568, "07"),
    @(23, 1300, 43.09749603271484, 4.174896240234375, 29.86320304870605, 57.64503860473633, 18438, "07"),
    @(24, 1400, 42.9649543762207,  4.170501232147217, 29.96230316162109, 56.96160125732422, 18366, "07"),
    @(25, 1500, 42.86069488525391, 4.255523204803467, 29.83586692810059, 60.26945114135742, 18392, "07")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).NumberFormat = "@"
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 7).ClearFormats()
}
